# Language.xlsx restructuring:
#  - Sheet1 "Sheet1" -> renamed "Comm", plus 4 new sheets (Property, Guild, Tip, Item)
#  - New/expanded shared strings for per-sheet UI text
#  - New cell content + sizing/selection per sheet
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheets: rename the original sheet, then create the rest in the exact
#    order that reproduces the target sheetId allocation (Comm=1, Property=7,
#    Guild=3, Tip=4, Item=5) and then move Property into final position.
# ---------------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item(1)
$wsComm.Name = "Comm"

$tmpA = $wb.Worksheets.Add($null, $wsComm)          # consumes sheetId 2 (discarded)
$wsGuild = $wb.Worksheets.Add($null, $tmpA)         # sheetId 3
$wsGuild.Name = "Guild"
$wsTip = $wb.Worksheets.Add($null, $wsGuild)        # sheetId 4
$wsTip.Name = "Tip"
$wsItem = $wb.Worksheets.Add($null, $wsTip)         # sheetId 5
$wsItem.Name = "Item"
$tmpB = $wb.Worksheets.Add($null, $wsItem)          # consumes sheetId 6 (discarded)
$wsProperty = $wb.Worksheets.Add($null, $tmpB)      # sheetId 7
$wsProperty.Name = "Property"

$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet1").Delete()

$wb.Worksheets.Item("Property").Move($wb.Worksheets.Item("Guild"))

# Re-fetch references now that the collection order changed.
$wsComm = $wb.Worksheets.Item("Comm")
$wsProperty = $wb.Worksheets.Item("Property")
$wsGuild = $wb.Worksheets.Item("Guild")
$wsTip = $wb.Worksheets.Item("Tip")
$wsItem = $wb.Worksheets.Item("Item")

# ---------------------------------------------------------------------------
# 2) Comm sheet: header already correct; rewrite the ID column (A) and the
#    Chinese column (C) with the new Comm-specific text; B keeps Langage_1..6.
# ---------------------------------------------------------------------------
$wsComm.Range("A2").Value = "Langage_Comm_1"
$wsComm.Range("A3").Value = "Langage_Comm_2"
$wsComm.Range("A4").Value = "Langage_Comm_3"
$wsComm.Range("A5").Value = "Langage_Comm_4"
$wsComm.Range("A6").Value = "Langage_Comm_5"
$wsComm.Range("A7").Value = "Langage_Comm_6"

$wsComm.Range("C2").Value = "确认"
$wsComm.Range("C3").Value = "取消"
$wsComm.Range("C4").Value = "登录"
$wsComm.Range("C5").Value = "创建角色"
$wsComm.Range("C6").Value = "进入游戏"
# C7 keeps its original text (中文_6)

# Extend with 5 more blank-but-styled rows (8-12), matching row 2's format.
$wsComm.Range("A2:C2").Copy($wsComm.Range("A8:C8"))
$wsComm.Range("A8:C8").ClearContents()
$wsComm.Range("A8:C8").Copy($wsComm.Range("A9:C12"))

$wsComm.Columns.Item(1).ColumnWidth = 31.1
$wsComm.Columns.Item(2).ColumnWidth = 23.65
$wsComm.Columns.Item(3).ColumnWidth = 22.3

$wsComm.Range("C8").Select()

# ---------------------------------------------------------------------------
# 3) Property sheet: single column (A) of property labels, header row B/C
#    still populated with English/Chinese labels (copied from Comm row 1).
# ---------------------------------------------------------------------------
$wsComm.Range("A1:C1").Copy($wsProperty.Range("A1:C1"))

$wsProperty.Range("A2").Value = "Langage_HP"
$wsProperty.Range("A3").Value = "Langage_MAXHP"
$wsProperty.Range("A4").Value = "Langage_MP"
$wsProperty.Range("A5").Value = "Langage_MAXMP"
$wsProperty.Range("A6").Value = "Langage_VP"
$wsProperty.Range("A7").Value = "Langage_ATTACK"

$wsProperty.Range("A2").Copy($wsProperty.Range("A8"))
$wsProperty.Range("A8").ClearContents()
$wsProperty.Range("A8").Copy($wsProperty.Range("A9:A28"))

$wsProperty.Columns.Item(1).ColumnWidth = 50.6

$wsProperty.Rows.Item(1).Select()

# ---------------------------------------------------------------------------
# 4) Guild sheet: one data row of guild-confirmation strings + blank styled
#    rows (3-12, 16) plus one isolated extra styled cell at row 22.
# ---------------------------------------------------------------------------
$wsComm.Range("A1:C1").Copy($wsGuild.Range("A1:C1"))
$wsComm.Range("A2:C2").Copy($wsGuild.Range("A2:C2"))

$wsGuild.Range("A2").Value = "Langage_Guild_1"
$wsGuild.Range("B2").Value = "Langage_1"
$wsGuild.Range("C2").Value = "确认要加入这个公会吗？点击确认加入"

$wsGuild.Range("A2:C2").Copy($wsGuild.Range("A3:C3"))
$wsGuild.Range("A3:C3").ClearContents()
$wsGuild.Range("A3:C3").Copy($wsGuild.Range("A4:C12"))
$wsGuild.Range("A3:C3").Copy($wsGuild.Range("A16:C16"))

$wsProperty.Range("A8").Copy($wsGuild.Range("A22"))

$wsGuild.Columns.Item(1).ColumnWidth = 31.1
$wsGuild.Columns.Item(2).ColumnWidth = 23.65
$wsGuild.Columns.Item(3).ColumnWidth = 22.3

$wsGuild.Range("A12").Select()

# ---------------------------------------------------------------------------
# 5) Tip / Item sheets: header row only.
# ---------------------------------------------------------------------------
$wsComm.Range("A1:C1").Copy($wsTip.Range("A1:C1"))
$wsTip.Rows.Item(1).Select()

$wsComm.Range("A1:C1").Copy($wsItem.Range("A1:C1"))
$wsItem.Rows.Item(1).Select()

# ---------------------------------------------------------------------------
# 6) Final UI state: Comm is the active/visible tab.
# ---------------------------------------------------------------------------
$wsComm.Activate()
$wsComm.Range("C8").Select()
